# Refresh the "cryptos" price/volume table (GitHub Actions scrape update).
# Price (D) / Volume(1h) (E) are stored as plain text in the sheet. Values
# that look like a bare number (e.g. "242.43") would otherwise be
# auto-coerced to a numeric cell by Excel, so those are entered with a
# leading apostrophe (forces text) and the cell style is then reset back
# to "Normal" so no stray quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.032.19'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.825.82'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''242.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '''0.6205'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").Value = '''0.9994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '''0.07406'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.12%  '
$ws.Range("D9").Value = '''0.2913'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '''23.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.35%  '
$ws.Range("D11").Value = '''0.07683'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = '1.829.63'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '''4.991'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").Value = '''0.6652'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '''82.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '''0.000009347'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.80%  '
$ws.Range("D17").Value = '''5.929'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("D18").Value = '29.031.40'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '2.074.93'
$ws.Range("D20").Value = '''12.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("D21").Value = '''222.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '''0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '''7.095'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("D24").Value = '''0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '''159.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").Value = '''0.1386'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("D27").Value = '''8.463'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = '''17.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").Value = '''1.486'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").Value = '''0.05718'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.65%  '
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '''4.118'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("E33").Value = '  +1.12%  '
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").Value = '''0.7366'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("D38").Value = '''2.757'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").Value = '1.217.16'
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").Value = '''0.01764'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").Value = '''6.467'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").Value = '''0.8901'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").Value = '''0.9990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '''101.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = '1.980.78'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").Value = '''0.00000000125'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.92%  '
$ws.Range("D47").Value = '''65.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").Value = '''0.07559'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +16.02%  '
$ws.Range("D50").Value = '''0.4033'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '''8.948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.08%  '
